# Slide 10: remove the "STRATEGY HUMAN RESOURCES:" textbox and the
# "Lorem ipsum..." rectangle that sit inside the "Group 92" group shape,
# then let the group's bounding box shrink to fit the remaining members.
#
# The runtime can't Delete() a shape addressed through GroupItems
# directly, so the group is ungrouped first, the two unwanted shapes are
# removed at the top level (where Delete works), and the surviving
# members are regrouped so the presentation keeps the same visual
# grouping it started with.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# Locate the group shape by name.
$grpIndex = -1
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
  if ($s.Shapes.Item($i).Name -eq "Group 92") {
    $grpIndex = $i
  }
}
$grp = $s.Shapes.Item($grpIndex)

# Ungroup so the nested shapes become directly addressable/deletable.
$ungrouped = $grp.Ungroup()

# Remove the two shapes called out in the edit: "TextBox 105" (the
# "STRATEGY HUMAN RESOURCES:" caption) and "Rectangle 106" (the Lorem
# ipsum paragraph).
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
  $nm = $s.Shapes.Item($i).Name
  if ($nm -eq "TextBox 105" -or $nm -eq "Rectangle 106") {
    $s.Shapes.Item($i).Delete()
  }
}

# Regroup the remaining former members ("TextBox 100", "Rectangle 102",
# "Straight Connector 107") back into a single group, restoring the
# original group's name. The group's extent now shrinks to wrap only
# the shapes that are left.
$names = @("TextBox 100", "Rectangle 102", "Straight Connector 107")
$newGrp = $s.Shapes.Range($names).Group()
$newGrp.Name = "Group 92"
